# Disaggregation of commodity Copper
#
# 1. Rename the shared commodity label "Copper ores and concentrates" -> "Copper"
#    on every year sheet (cell C4, which is merged B-column category / C-column
#    commodity name, repeated identically across all "20xx" sheets).
# 2. A handful of sheets have their D4 (Copper total) value refreshed to a
#    slightly different floating point figure after the disaggregation.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C4").Value = "Copper"
}

$valueUpdates = @(
    @("2033", 95388.41488819558),
    @("2047", 634503.609349301),
    @("2048", 806653.9762728701),
    @("2054", 1998651.653451595),
    @("2065", 909749.2385804425),
    @("2073", 879339.2488812557)
)

foreach ($pair in $valueUpdates) {
    $sheet = $wb.Worksheets.Item($pair[0])
    $sheet.Range("D4").Value = $pair[1]
}
